$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "Datos actualizados a 14 de Abril de 2020 a las 08:22" -> "Datos actualizados a 14 de Abril de 2020 a las 08:52"
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Abril de 2020 a las 08:52"

# Row 21: "Israel" (data update)
$ws.Cells.Item(21, 2).Value = 11868
$ws.Cells.Item(21, 3).Value = 282
$ws.Cells.Item(21, 4).Value = 2000
$ws.Cells.Item(21, 5).Value = 9751
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 117

# Row 35: "Chequia" (data update)
$ws.Cells.Item(35, 4).Value = 527
$ws.Cells.Item(35, 5).Value = 5385
$ws.Cells.Item(35, 6).Value = 92
$ws.Cells.Item(35, 7).Value = 4
$ws.Cells.Item(35, 8).Value = 147

# Row 45: "Luxemburgo" -> "Ucrania"
$ws.Cells.Item(45, 1).Value = "Ucrania"
$ws.Cells.Item(45, 2).Value = 3372
$ws.Cells.Item(45, 3).Value = 270
$ws.Cells.Item(45, 4).Value = 119
$ws.Cells.Item(45, 5).Value = 3155
$ws.Cells.Item(45, 6).Value = 45
$ws.Cells.Item(45, 7).Value = 5
$ws.Cells.Item(45, 8).Value = 98

# Row 46: "Catar" -> "Luxemburgo"
$ws.Cells.Item(46, 1).Value = "Luxemburgo"
$ws.Cells.Item(46, 2).Value = 3292
$ws.Cells.Item(46, 4).Value = 500
$ws.Cells.Item(46, 5).Value = 2723
$ws.Cells.Item(46, 6).Value = 30
$ws.Cells.Item(46, 8).Value = 69

# Row 47: "Republica Dominicana" -> "Catar"
$ws.Cells.Item(47, 1).Value = "Catar"
$ws.Cells.Item(47, 2).Value = 3231
$ws.Cells.Item(47, 4).Value = 334
$ws.Cells.Item(47, 5).Value = 2890
$ws.Cells.Item(47, 6).Value = 37
$ws.Cells.Item(47, 8).Value = 7

# Row 48: "Ucrania" -> "Republica Dominicana"
$ws.Cells.Item(48, 1).Value = "Republica Dominicana"
$ws.Cells.Item(48, 2).Value = 3167
$ws.Cells.Item(48, 4).Value = 152
$ws.Cells.Item(48, 5).Value = 2838
$ws.Cells.Item(48, 6).Value = 147
$ws.Cells.Item(48, 8).Value = 177

# Row 79: "Banglades" -> "Oman"
$ws.Cells.Item(79, 1).Value = "Oman"
$ws.Cells.Item(79, 2).Value = 813
$ws.Cells.Item(79, 3).Value = 86
$ws.Cells.Item(79, 4).Value = 130
$ws.Cells.Item(79, 5).Value = 679
$ws.Cells.Item(79, 6).Value = 3
$ws.Cells.Item(79, 8).Value = 4

# Row 80: "Eslovaquia" -> "Banglades"
$ws.Cells.Item(80, 1).Value = "Banglades"
$ws.Cells.Item(80, 2).Value = 803
$ws.Cells.Item(80, 4).Value = 42
$ws.Cells.Item(80, 5).Value = 722
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 8).Value = 39

# Row 81: "Oman" -> "Eslovaquia"
$ws.Cells.Item(81, 1).Value = "Eslovaquia"
$ws.Cells.Item(81, 2).Value = 769
$ws.Cells.Item(81, 4).Value = 107
$ws.Cells.Item(81, 5).Value = 660
$ws.Cells.Item(81, 6).Value = 5
$ws.Cells.Item(81, 8).Value = 2

# Row 111: "Senegal" -> "Georgia"
$ws.Cells.Item(111, 1).Value = "Georgia"
$ws.Cells.Item(111, 2).Value = 296
$ws.Cells.Item(111, 3).Value = 24
$ws.Cells.Item(111, 4).Value = 68
$ws.Cells.Item(111, 5).Value = 225
$ws.Cells.Item(111, 6).Value = 6
$ws.Cells.Item(111, 8).Value = 3

# Row 112: "Montenegro" -> "Senegal"
$ws.Cells.Item(112, 1).Value = "Senegal"
$ws.Cells.Item(112, 2).Value = 291
$ws.Cells.Item(112, 4).Value = 178
$ws.Cells.Item(112, 5).Value = 111
$ws.Cells.Item(112, 6).Value = 1
$ws.Cells.Item(112, 8).Value = 2

# Row 113: "Georgia" -> "Montenegro"
$ws.Cells.Item(113, 1).Value = "Montenegro"
$ws.Cells.Item(113, 2).Value = 274
$ws.Cells.Item(113, 4).Value = 5
$ws.Cells.Item(113, 5).Value = 266
$ws.Cells.Item(113, 6).Value = 7

# Row 121: "Islas Feroe" (data update)
$ws.Cells.Item(121, 4).Value = 163
$ws.Cells.Item(121, 5).Value = 21

# Row 204: "San Bartolome" -> "Timor Oriental"
$ws.Cells.Item(204, 1).Value = "Timor Oriental"
$ws.Cells.Item(204, 3).Value = 2
$ws.Cells.Item(204, 4).Value = 1
$ws.Cells.Item(204, 5).Value = 5

# Row 205: "Islas Malvinas" -> "San Bartolome"
$ws.Cells.Item(205, 1).Value = "San Bartolome"
$ws.Cells.Item(205, 2).Value = 6
$ws.Cells.Item(205, 4).Value = 4
$ws.Cells.Item(205, 5).Value = 2

# Row 207: "Butan" -> "Islas Malvinas"
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"
$ws.Cells.Item(207, 4).Value = 1
$ws.Cells.Item(207, 5).Value = 4

# Row 208: "Sudan del Sur" -> "Butan"
$ws.Cells.Item(208, 1).Value = "Butan"
$ws.Cells.Item(208, 2).Value = 5
$ws.Cells.Item(208, 4).Value = 2
$ws.Cells.Item(208, 5).Value = 3

# Row 210: "Timor Oriental" -> "Sudan del Sur"
$ws.Cells.Item(210, 1).Value = "Sudan del Sur"
$ws.Cells.Item(210, 4).Value = 0
$ws.Cells.Item(210, 5).Value = 4
